$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.847.27"
$ws.Range("E2").Value = "  +1.59%  "
$ws.Range("D3").Value = "1.725.93"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9976"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9982"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4885"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.52%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2589"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06208"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").Value = "1.732.73"
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "15.98"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06901"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6069"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.474"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.20"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9977"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "26.627.05"
$ws.Range("E17").Value = "  +0.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9977"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007175"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("E20").Value = "  +0.85%  "
$ws.Range("D21").Value = "1.950.54"
$ws.Range("E21").Value = "  +0.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.413"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.575"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.082"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.771"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "106.33"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.63%  "
$ws.Range("E29").Value = "  -1.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.947"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07982"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.683"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04519"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.28%  "
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.008"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6247"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9337"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.042"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.445"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9975"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01497"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.73%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.671"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3840"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.847"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1161"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05399"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.902"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.70%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.19"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "51.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.233"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.27%  "
